# GreatLink Asia High Dividend Equity - DividendHistory update
# Adds the newest dividend record (XD Date 27/03/2025) as a new row right
# below the header, pushing all the historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (just below the header row), shifting
# all existing dividend history rows down by one.
$ws.Rows("2:2").Insert()

# Row 3 now holds what used to be row 2 (30/12/2024 | 30/12/2024 | 0.008).
# Copy its Gross Dividend cell down into the new row first so the new
# C2 cell inherits the plain "text number" formatting already used
# throughout the column (avoids Excel re-interpreting "0.008" as a
# numeric value when it is typed in directly).
$ws.Range("C3").Copy($ws.Range("C2"))

# Fill in the new dividend record.
$ws.Range("A2").Value = "27/03/2025"
$ws.Range("B2").Value = "27/03/2025"
